$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.501.33"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.121.28"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.39"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.56"
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.118.57"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.442"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").Value = "3.653.77"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.135"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.34"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "57.613.00"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "3.113.33"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.01"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "352.15"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.30"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.506"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "0.0₃0917"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.25"
$ws.Range("E31").Value = "  -6.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.12"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.87"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.25"
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.13"
$ws.Range("E38").Value = "  -6.13%  "
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0666"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.61"
$ws.Range("E41").Value = "  +5.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.14"
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.697"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "3.159.14"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.34"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0267"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "2.318.21"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.971"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.04"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.15"
$ws.Range("E51").Value = "  -4.79%  "
